$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update "Förändrad" (column C) for rows 2-25 from 46073 to 46074
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 3).Value = 46074
}

# 2. Rows 7-25 got reshuffled: each row's (A, B, G) data - and for one row, the F
#    "Kyrkan" marker - moved to a different row. Capture the current A/B/G/F
#    values for rows 7-25 first, then write them back out in the new order.

$srcRows = @(17, 21, 7, 8, 25, 19, 23, 11, 10, 22, 24, 13, 14, 15, 20, 12, 18, 16, 9)

$dataA = @{}
$dataB = @{}
$dataG = @{}
$dataF = @{}
foreach ($r in 7..25) {
    $dataA[$r] = $ws.Cells.Item($r, 1).Value()
    $dataB[$r] = $ws.Cells.Item($r, 2).Value()
    $dataG[$r] = $ws.Cells.Item($r, 7).Value()
    $dataF[$r] = $ws.Cells.Item($r, 6).Value()
}

$destRows = 7..25
for ($i = 0; $i -lt $destRows.Length; $i++) {
    $dest = $destRows[$i]
    $src = $srcRows[$i]
    $ws.Cells.Item($dest, 1).Value = $dataA[$src]
    $ws.Cells.Item($dest, 2).Value = $dataB[$src]
    $ws.Cells.Item($dest, 7).Value = $dataG[$src]

    $fVal = $dataF[$src]
    if ($fVal -ne $null -and $fVal -ne "") {
        $ws.Cells.Item($dest, 6).Value = $fVal
    } else {
        $ws.Cells.Item($dest, 6).Value = ""
    }
}
